$d = $word.ActiveDocument

$replacements = @(
    @{old="238×8="; new="747×3="},
    @{old="234×4="; new="789×9="},
    @{old="979×9="; new="782×5="},
    @{old="456×5="; new="344×2="},
    @{old="182×3="; new="172×2="},
    @{old="129×6="; new="404×6="},
    @{old="343×6="; new="191×4="},
    @{old="356×2="; new="406×5="},
    @{old="280×3="; new="788×5="},
    @{old="935×5="; new="275×6="},
    @{old="174×3="; new="135×2="},
    @{old="345×6="; new="531×7="},
    @{old="144×8="; new="782×5="},
    @{old="167×2="; new="764×9="},
    @{old="234×7="; new="572×8="},
    @{old="120×2="; new="534×9="},
    @{old="233×8="; new="509×5="},
    @{old="737×3="; new="756×2="},
    @{old="292×8="; new="195×3="},
    @{old="803×9="; new="912×7="},
    @{old="262×2="; new="550×6="},
    @{old="990×4="; new="677×2="},
    @{old="523×5="; new="543×3="},
    @{old="214×3="; new="943×4="},
    @{old="785×4="; new="574×4="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
